# Estado de cuenta NIT-9000832838: se eliminan los periodos de mora
# anteriores y se agregan nuevos; se actualiza la base de datos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: row 35 (old) carried the "closing border" style because it
#    used to be the last data row. We now have one extra data row, so that
#    bordered style must move down to the new last row (36), and row 35
#    becomes a normal interior row (copy row 34's look).
# ---------------------------------------------------------------------------
$ws.Range("B35:J35").Copy()
$ws.Range("B36:J36").PasteSpecial(-4122)

$ws.Range("B34:J34").Copy()
$ws.Range("B35:J35").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Rewrite the worker/mora table, rows 16-36.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora,
#             G=Salario Basico
# ---------------------------------------------------------------------------
$data = @(
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2110",29260,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2109",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2108",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2107",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2106",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2105",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2104",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2103",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2102",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2101",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2012",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2011",35112,877803),
  @("CC","73163175","JESUS MARIA MEZA JACKSON","2010",19897,877803),
  @("CC","9271328","ALVARO RUIZ OLIVEROS","2102",35112,908526),
  @("CC","73183653","GABRIEL ENRIQUE VASQUEZ REYES","2102",35112,908526),
  @("CC","1047460323","PEDRO LUIS VASQUEZ REYES","2102",35112,908526),
  @("CC","1143401546","YOVANIS MORALES CASTELLAR","2102",35112,1423500),
  @("CC","1047495328","JESUS RAFAEL ANILLO RIOS","2102",35112,908526),
  @("CC","1047488768","GERMAIN BARRAGAN CAMPILLO","2102",35112,908526),
  @("CC","107389299","GUSTAVO ADOLFO VARGAS GONZALEZ","2011",35112,877803),
  @("CC","107389299","GUSTAVO ADOLFO VARGAS GONZALEZ","2010",19897,877803)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r, 2).Value2 = $row[0]
  $ws.Cells.Item($r, 3).Value2 = $row[1]
  $ws.Cells.Item($r, 4).Value2 = $row[2]
  $ws.Cells.Item($r, 5).Value2 = $row[3]
  $ws.Cells.Item($r, 6).Value2 = $row[4]
  $ws.Cells.Item($r, 7).Value2 = $row[5]
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Summary header fields.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 701070   # VALOR MORA total
$ws.Range("C13").Value2 = 8        # Cant. Trabajadores

# ---------------------------------------------------------------------------
# 4) Shift the signature block down one row:
#    old row 40 (line) / row 41 (labels)  ->  row 41 (line) / row 42 (labels)
# ---------------------------------------------------------------------------
$ws.Range("B40:C40").UnMerge()
$ws.Range("H40:J40").UnMerge()
$ws.Range("B41:C41").UnMerge()
$ws.Range("H41:J41").UnMerge()

$ws.Range("B41:C41").Copy()
$ws.Range("B42:C42").PasteSpecial(-4122)
$ws.Range("H41:J41").Copy()
$ws.Range("H42:J42").PasteSpecial(-4122)

$ws.Range("B40:C40").Copy()
$ws.Range("B41:C41").PasteSpecial(-4122)
$ws.Range("H40:J40").Copy()
$ws.Range("H41:J41").PasteSpecial(-4122)

$ws.Range("B40:J40").Clear()

$ws.Range("B41:C41").Merge()
$ws.Range("H41:J41").Merge()
$ws.Range("B42:C42").Merge()
$ws.Range("H42:J42").Merge()

$ws.Range("B41").Value2 = "___________________________________"
$ws.Range("H41").Value2 = "___________________________________"
$ws.Range("B42").Value2 = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H42").Value2 = "FIRMA DEL REPRESENTANTE LEGAL"

Write-Output "edit applied"
